$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The shared-string table was regenerated upstream so that a new "Holden"
# series was inserted ahead of the existing "HexGrid-90degTilt*degRes"
# series. Because the worksheet rows still point at the same (now shifted)
# string-table slots, rows 16-19 end up labelled with the new Holden
# series, and the displaced HexGrid labels reappear on four brand new rows
# (20-23) appended at the bottom with the same "all ones" data pattern.

# 1) Relabel rows 16-19 from the HexGrid series to the new Holden series.
$ws.Range("B16").Value = "Holden2.5"
$ws.Range("B17").Value = "Holden5"
$ws.Range("B18").Value = "Holden10"
$ws.Range("B19").Value = "Holden15"

# 2) Drop the duplicated tail columns (X:AG) on rows 1-2; the sheet's used
#    range shrinks from A1:AG19 to A1:W23.
$ws.Range("X1:AG2").Clear()

# 3) Append the four displaced HexGrid rows at the bottom of the table.
$newRows = @(
    @{ Row = 20; Id = 18; Name = "HexGrid-90degTilt2.5degRes" },
    @{ Row = 21; Id = 19; Name = "HexGrid-90degTilt5degRes" },
    @{ Row = 22; Id = 20; Name = "HexGrid-90degTilt10degRes" },
    @{ Row = 23; Id = 21; Name = "HexGrid-90degTilt15degRes" }
)

foreach ($r in $newRows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value = $r.Id
    $ws.Cells.Item($row, 2).Value = $r.Name

    for ($col = 3; $col -le 23; $col++) {
        $ws.Cells.Item($row, $col).Value = 1
    }

    # Match the formatting already used by the other data rows: column A
    # is bold/centered/bordered (style index 1), column B..W plain.
    $ws.Range("A19").Copy()
    $ws.Range("A$row").PasteSpecial(-4122)
}

$excel.CutCopyMode = 0
$wb.Save()
